# [INITIAL COMMIT] - Filter datatable within number of X days, and also
# added handling for rows with empty birthday.
#
# Updates the "Settings" (BirthdayList_*) rows and adds the new
# "UpcomingBirthdayRange" setting plus a new "Exception_BirthdayListNotExist"
# constant on the "Constants" sheet. Also moves the active tab / selection
# from Settings!C9 to Constants!C18 (with Settings left selected at B10).

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Settings")
$ws2 = $wb.Worksheets.Item("Constants")

# --- Settings sheet -------------------------------------------------------

# Row 6 (BirthdayList_Path): reword the description and update the sample
# path to use the {{UserProfile}} placeholder; the value cell now wraps.
$ws1.Range("C6").Value = "Filepath of the excel spreadsheet. Where UserProfile is the name of the logged user"
$ws1.Range("B6").Value = "C:\Users\{{UserProfile}}\Documents\UiPath\Birthday Greeting Postcard\C-ET Birthday List.xlsx"
$ws1.Range("B6").WrapText = $true

# Row 8 (BirthdayList_Range): clear the literal "" default value and reword
# the description.
$ws1.Range("B8").ClearContents()
$ws1.Range("C8").Value = "Range to be read in the birthday list excel spreadsheet. Keep the value emptied to read all cells in the spreadsheet"

# Row 10 (new setting): UpcomingBirthdayRange = 14
$ws1.Range("A10").Value = "UpcomingBirthdayRange"
$ws1.Range("B10").Value = 14
$ws1.Range("C10").Value = "Value in number of days to determine upcoming birthdays within range (DateToday - X days)"

# --- Constants sheet --------------------------------------------------------

# Row 12 (new constant): Exception_BirthdayListNotExist
$ws2.Range("A12").Value = "Exception_BirthdayListNotExist"
$ws2.Range("C12").Value = "Exception message if the list of birthday excel file does not exist"
$ws2.Range("B12").Value = "Birthday master list does not exist"

# --- Active sheet / selection ---------------------------------------------

$ws1.Range("B10").Select()
$ws2.Select()
$ws2.Range("C18").Select()
